$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $pattern) {
    $idx = 0
    $found = -1
    foreach ($pp in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($pp.Range.Text -like $pattern) {
            $found = $idx
        }
    }
    return $found
}

# ---------------------------------------------------------------------------
# 1. These four items in the top "to-do" list are now finished (the logic
#    moved into firebase.tsx / risk cards got fully implemented), so remove
#    them from here - they get re-added to the "Done:" list in step 4.
# ---------------------------------------------------------------------------
$completedItems = @(
    "Lose a turn functionality*",
    "No rent due functionality*",
    "Multiple players on same space?*",
    "Risk Cards*"
)
foreach ($pattern in $completedItems) {
    $idx = Get-ParagraphIndexByText $d $pattern
    if ($idx -gt 0) {
        $d.Paragraphs($idx).Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. "Trade property with other players" functionality was moved into
#    firebase.tsx, so the to-do entry is removed entirely (not re-added
#    anywhere).
# ---------------------------------------------------------------------------
$tradeIndex = Get-ParagraphIndexByText $d "*Trade property with other players*"
if ($tradeIndex -gt 0) {
    $d.Paragraphs($tradeIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Tidy up the "Pass go = 50*properties owned" entry in the "Done:" list -
#    drop the trailing space and the stray tab character after it.
# ---------------------------------------------------------------------------
[void]$d.Content.Find.Execute("owned " + [char]9, $true, $false, $false, $false, $false, `
    $true, 1, $false, "owned", 2)

# ---------------------------------------------------------------------------
# 4. Mark the four completed items from step 1 as done by appending them
#    (with strikethrough formatting) to the end of the "Done:" list, right
#    after "Pass go = 50*properties owned".
# ---------------------------------------------------------------------------
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$pPrStrike = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr>'

$doneParagraphsXml = @(
    ('<w:p>' + $pPrStrike + '<w:r><w:rPr><w:strike/></w:rPr><w:t>Lose a turn functionality</w:t></w:r></w:p>'),
    ('<w:p>' + $pPrStrike + '<w:r><w:rPr><w:strike/></w:rPr><w:t>No rent due functionality</w:t></w:r></w:p>'),
    ('<w:p>' + $pPrStrike + '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Multiple players </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>on same</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> space?</w:t></w:r></w:p>'),
    ('<w:p>' + $pPrStrike + '<w:r><w:rPr><w:strike/></w:rPr><w:t>Risk Cards</w:t></w:r></w:p>')
)

$anchorIndex = Get-ParagraphIndexByText $d "Pass go = 50*properties owned*"
foreach ($pXml in $doneParagraphsXml) {
    $d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()
    $anchorIndex = $anchorIndex + 1
    [void]$d.Paragraphs($anchorIndex).Range.InsertXML($xmlHeader + $pXml + $xmlFooter)
}
